# Update column G (daily value) for rows 2-6 and recompute the AG total
# column (sum of B:AF) to match the new figures reported by BIBI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newG = @{
    2 = 5403.09
    3 = 8020.9
    4 = 1287
    5 = 1823
    6 = 16533.99
}

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 7).Value = $newG[$row]   # Column G is the 7th column

    # Recalculate the row total (column AG, 33rd column) as the sum of B:AF
    $total = 0
    for ($col = 2; $col -le 32; $col++) {
        $total += $ws.Cells.Item($row, $col).Value2
    }
    $ws.Cells.Item($row, 33).Value = $total
}
